# Weekly update: a new price-report row for Cereza (Agrícola del Norte S.A. de Arica)
# is inserted at the top of the data block (row 19), pushing the existing
# rows 19-31 down to 20-32 (the sheet's dimension grows from T31 to T32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 19 - this shifts rows 19:31 down to 20:32
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with this week's record
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 45280
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103001
$ws.Range("J19").Value = "Cereza"
$ws.Range("K19").Value = "Santina"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 12000
$ws.Range("P19").Value = 11000
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1100
$ws.Range("T19").Value = 10
